$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet and update header label to reflect new "through" date
$ws.Name = "Through 2022-05-26"
$ws.Range("I1").Value = "2022 (through 05-26)"

# Update the June total (row 6) and the grand total (row 14)
$ws.Range("I6").Value = 95
$ws.Range("I14").Value = 646
